$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text code that looks numeric ("001"); force text storage,
# then clear the formatting override so no stray style is left behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# M2 / N2 are plain text date-time strings
$ws.Range("M2").Value = "2020-12-22 00:00:00"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

# Numeric financial fields
$ws.Range("O2").Value = 388395121.71
$ws.Range("P2").Value = 87021828.09999999
$ws.Range("Q2").Value = 139974760.49
$ws.Range("R2").Value = -2.7549123139
$ws.Range("S2").Value = 48411989.22
$ws.Range("T2").Value = 90.5041797141
$ws.Range("U2").Value = 35137931.17
$ws.Range("V2").Value = -17.2065907378
$ws.Range("W2").Value = 63777399
$ws.Range("X2").Value = 31559983.37
$ws.Range("Y2").Value = 118.5227290454
$ws.Range("Z2").Value = 4370015.36
$ws.Range("AA2").Value = -58.545184113
$ws.Range("AB2").Value = 324617722.71
$ws.Range("AC2").Value = 3.4470545833
$ws.Range("AD2").Value = 4.1272129932
$ws.Range("AE2").Value = 7.7325365924
$ws.Range("AF2").Value = 473.609706724
$ws.Range("AG2").Value = 16.4207518156
